# "added requested features" — add a new H column (SEALED flag) to the
# product table and replace the numeric city_ID (column D) with the
# actual product/order code strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H
$ws.Range("H4").Value = "SEALED"

# Row -> (new product code for column D, new SEALED flag for column H or $null to leave blank)
$rowData = @{
     5 = @{ Code = "K1L203"; Sealed = $true  }
     6 = @{ Code = "K2P032"; Sealed = $null  }
     7 = @{ Code = "Q20F55"; Sealed = $false }
     8 = @{ Code = "K2P032"; Sealed = $null  }
     9 = @{ Code = "Q20F55"; Sealed = $true  }
    10 = @{ Code = "K1L203"; Sealed = $null  }
    11 = @{ Code = "UI032P"; Sealed = $false }
    12 = @{ Code = "UI032P"; Sealed = $null  }
    13 = @{ Code = "K2P032"; Sealed = $true  }
    14 = @{ Code = "UI032P"; Sealed = $null  }
    15 = @{ Code = "K1L203"; Sealed = $false }
    16 = @{ Code = "Q20F55"; Sealed = $null  }
}

foreach ($r in 5..16) {
    $info = $rowData[$r]
    $ws.Cells.Item($r, 4).Value = $info.Code
    if ($null -ne $info.Sealed) {
        $ws.Cells.Item($r, 8).Value = $info.Sealed
    }
}

# Matches the author's final selection in the saved workbook.
$ws.Range("B9").Select() | Out-Null
